$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.67
$ws.Range("I2").Value = 5.5
$ws.Range("J2").Value = 2.3
$ws.Range("K2").Value = 2.2
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 3.4
$ws.Range("Q2").Value = 2.05
$ws.Range("R2").Value = 1.75
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 1.72
$ws.Range("Z2").Value = 12
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 7
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 67
$ws.Range("AU2").Value = 9

# Row 3
$ws.Range("G3").Value = 1.85
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 2.6
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("AN3").Value = 3.6
$ws.Range("AU3").Value = 9.5
$ws.Range("BB3").Value = 401

# Row 5
$ws.Range("G5").Value = 2.88
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 1.27
$ws.Range("X5").Value = 15
$ws.Range("Y5").Value = 11
$ws.Range("AN5").Value = 5
$ws.Range("AO5").Value = 17

# Row 6
$ws.Range("P6").Value = 3.92
$ws.Range("U6").Value = 1.88
$ws.Range("V6").Value = 1.88

# Row 7
$ws.Range("G7").Value = 2.8
$ws.Range("H7").Value = 3.1
$ws.Range("I7").Value = 2.45
$ws.Range("J7").Value = 3.25
$ws.Range("L7").Value = 3.05
$ws.Range("M7").Value = 10.4
$ws.Range("Q7").Value = 1.78
$ws.Range("R7").Value = 1.93
$ws.Range("U7").Value = 1.55
$ws.Range("V7").Value = 2.15
$ws.Range("W7").Value = 11
$ws.Range("X7").Value = 17
$ws.Range("Y7").Value = 9.75
$ws.Range("Z7").Value = 37
$ws.Range("AA7").Value = 22
$ws.Range("AB7").Value = 24
$ws.Range("AE7").Value = 11.25
$ws.Range("AH7").Value = 9
$ws.Range("AI7").Value = 13
$ws.Range("AJ7").Value = 9.25
$ws.Range("AK7").Value = 28
$ws.Range("AL7").Value = 19.5
$ws.Range("AN7").Value = 4.85
$ws.Range("AO7").Value = 14.5
$ws.Range("AP7").Value = 19
$ws.Range("AQ7").Value = 65
$ws.Range("AR7").Value = 80
$ws.Range("AS7").Value = 200
$ws.Range("AW7").Value = 4.5
$ws.Range("AX7").Value = 13.5
$ws.Range("AY7").Value = 19
$ws.Range("AZ7").Value = 55
$ws.Range("BA7").Value = 80

# Row 8
$ws.Range("M8").Value = 1.03
$ws.Range("O8").Value = 1.22
$ws.Range("U8").Value = 1.73
$ws.Range("BD8").Value = 126

# Row 9
$ws.Range("G9").Value = 1.53
$ws.Range("H9").Value = 4.1
$ws.Range("J9").Value = 2.05
$ws.Range("K9").Value = 2.4
$ws.Range("M9").Value = 1.02
$ws.Range("N9").Value = 15
$ws.Range("O9").Value = 1.15
$ws.Range("P9").Value = 4.5
$ws.Range("Q9").Value = 1.62
$ws.Range("R9").Value = 2.25
$ws.Range("S9").Value = 1.3
$ws.Range("T9").Value = 3.4
$ws.Range("U9").Value = 1.73
$ws.Range("V9").Value = 2
$ws.Range("W9").Value = 8.5
$ws.Range("X9").Value = 8
$ws.Range("AB9").Value = 21
$ws.Range("AC9").Value = 15
$ws.Range("AD9").Value = 8
$ws.Range("AE9").Value = 15
$ws.Range("AF9").Value = 41
$ws.Range("AG9").Value = 151
$ws.Range("AH9").Value = 19
$ws.Range("AK9").Value = 51
$ws.Range("AP9").Value = 17
$ws.Range("AQ9").Value = 21
$ws.Range("AS9").Value = 101
$ws.Range("AT9").Value = 3.4
$ws.Range("AU9").Value = 8
$ws.Range("AW9").Value = 7.5
$ws.Range("AY9").Value = 29
$ws.Range("BA9").Value = 101

# Row 10
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 4.1
$ws.Range("J10").Value = 2.6
$ws.Range("K10").Value = 2.1
$ws.Range("M10").Value = 1.05
$ws.Range("N10").Value = 9
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = 2.1
$ws.Range("R10").Value = 1.7
$ws.Range("S10").Value = 1.44
$ws.Range("T10").Value = 2.63
$ws.Range("U10").Value = 1.91
$ws.Range("V10").Value = 1.8
$ws.Range("W10").Value = 6.5
$ws.Range("X10").Value = 8.5
$ws.Range("AB10").Value = 29
$ws.Range("AC10").Value = 8.5
$ws.Range("AE10").Value = 17
$ws.Range("AF10").Value = 51
$ws.Range("AL10").Value = 34
$ws.Range("AR10").Value = 51
$ws.Range("AT10").Value = 2.63
$ws.Range("AU10").Value = 8.5
$ws.Range("BA10").Value = 101
$ws.Range("BB10").Value = 251

# Row 11
$ws.Range("G11").Value = 2
$ws.Range("I11").Value = 3.8
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 7.5
$ws.Range("O11").Value = 1.41
$ws.Range("P11").Value = 2.62
$ws.Range("S11").Value = 1.53
$ws.Range("T11").Value = 2.38
$ws.Range("AH11").Value = 8.5
$ws.Range("AI11").Value = 17
$ws.Range("AO11").Value = 12
$ws.Range("AT11").Value = 2.38

# Row 12
$ws.Range("M12").Value = 1.05
$ws.Range("N12").Value = 8
$ws.Range("O12").Value = 1.41
$ws.Range("P12").Value = 2.62

# Row 13
$ws.Range("M13").Value = 1.03
$ws.Range("O13").Value = 1.22

Write-Host "Applied 148 cell updates"